$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The DIRECTION/INCLINATION pair was recorded with inclination measured from
# the horizontal using the opposite sign convention (welly/pandas bug).
# Re-derive INCLINATION (column D, rows 2:33) by adding 90 to each existing
# value so it matches the corrected convention.
for ($r = 2; $r -le 33; $r++) {
    $cell = $ws.Cells.Item($r, 4)
    $cell.Value2 = $cell.Value2 + 90
}

# Leave the sheet with the same selection state recorded after the fix.
$ws.Range("K14").Select()
